$ws = $excel.ActiveWorkbook.ActiveSheet

$ws.Range("G2").Value = 7.620274999999999
$ws.Range("H2").Value = 22.860825
$ws.Range("I2").Value = 0.6584612850834004
$ws.Range("J2").Value = 0.6584612850834003
$ws.Range("M2").Value = 10.055569
$ws.Range("N2").Value = 30.166707
$ws.Range("O2").Value = 0.349442268297237
$ws.Range("P2").Value = 0.3494422682972371
$ws.Range("Q2").Value = 76.626201061475
$ws.Range("R2").Value = 689.635809553275
$ws.Range("S2").Value = 0.2300942050454571
$ws.Range("T2").Value = 0.2300942050454571

$ws.Range("G3").Value = 7.620274999999999
$ws.Range("H3").Value = 22.860825
$ws.Range("I3").Value = 0.6584612850834004
$ws.Range("J3").Value = 0.6584612850834003
$ws.Range("O3").Value = 0.2478120087748427
$ws.Range("P3").Value = 0.2478120087748427
$ws.Range("Q3").Value = 54.34057219911666
$ws.Range("R3").Value = 489.06514979205
$ws.Range("S3").Value = 0.1631746137569818
$ws.Range("T3").Value = 0.1631746137569818

$ws.Range("G4").Value = 7.620274999999999
$ws.Range("H4").Value = 22.860825
$ws.Range("I4").Value = 0.6584612850834004
$ws.Range("J4").Value = 0.6584612850834003
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.623176333333333
$ws.Range("N4").Value = 7.869529
$ws.Range("O4").Value = 0.09115831118692827
$ws.Range("P4").Value = 0.09115831118692828
$ws.Range("Q4").Value = 19.98932503349166
$ws.Range("R4").Value = 179.903925301425
$ws.Range("S4").Value = 0.0600242187301773
$ws.Range("T4").Value = 0.0600242187301773

$ws.Range("G5").Value = 7.620274999999999
$ws.Range("H5").Value = 22.860825
$ws.Range("I5").Value = 0.6584612850834004
$ws.Range("J5").Value = 0.6584612850834003
$ws.Range("M5").Value = 6.063478666666666
$ws.Range("N5").Value = 18.190436
$ws.Range("O5").Value = 0.2107126646987263
$ws.Range("P5").Value = 0.2107126646987263
$ws.Range("Q5").Value = 46.20537489663332
$ws.Range("R5").Value = 415.8483740697
$ws.Range("S5").Value = 0.138746131980871
$ws.Range("T5").Value = 0.138746131980871

$ws.Range("G6").Value = 7.620274999999999
$ws.Range("H6").Value = 22.860825
$ws.Range("I6").Value = 0.6584612850834004
$ws.Range("J6").Value = 0.6584612850834003
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.5298106666666667
$ws.Range("N6").Value = 1.589432
$ws.Range("O6").Value = 0.01841151317524362
$ws.Range("P6").Value = 0.01841151317524363
$ws.Range("Q6").Value = 4.037302977933333
$ws.Range("R6").Value = 36.3357268014
$ws.Range("S6").Value = 0.01212326862570087
$ws.Range("T6").Value = 0.01212326862570087

$ws.Range("G7").Value = 7.620274999999999
$ws.Range("H7").Value = 22.860825
$ws.Range("I7").Value = 0.6584612850834004
$ws.Range("J7").Value = 0.6584612850834003
$ws.Range("M7").Value = 2.372966333333333
$ws.Range("N7").Value = 7.118898999999999
$ws.Range("O7").Value = 0.08246323386702208
$ws.Range("P7").Value = 0.08246323386702209
$ws.Range("Q7").Value = 18.08265602574166
$ws.Range("R7").Value = 162.743904231675
$ws.Range("S7").Value = 0.05429884694421235
$ws.Range("T7").Value = 0.05429884694421235

$ws.Range("I8").Value = 0.262323813236933
$ws.Range("J8").Value = 0.262323813236933
$ws.Range("M8").Value = 10.055569
$ws.Range("N8").Value = 30.166707
$ws.Range("O8").Value = 0.349442268297237
$ws.Range("P8").Value = 0.3494422682972371
$ws.Range("Q8").Value = 30.52704496325867
$ws.Range("R8").Value = 274.743404669328
$ws.Range("S8").Value = 0.09166702832589466
$ws.Range("T8").Value = 0.09166702832589467

$ws.Range("I9").Value = 0.262323813236933
$ws.Range("J9").Value = 0.262323813236933
$ws.Range("O9").Value = 0.2478120087748427
$ws.Range("P9").Value = 0.2478120087748427
$ws.Range("S9").Value = 0.06500699110772105
$ws.Range("T9").Value = 0.06500699110772107

$ws.Range("I10").Value = 0.262323813236933
$ws.Range("J10").Value = 0.262323813236933
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.623176333333333
$ws.Range("N10").Value = 7.869529
$ws.Range("O10").Value = 0.09115831118692827
$ws.Range("P10").Value = 0.09115831118692828
$ws.Range("Q10").Value = 7.963529649512889
$ws.Range("R10").Value = 71.67176684561601
$ws.Range("S10").Value = 0.023912995798794
$ws.Range("T10").Value = 0.023912995798794

$ws.Range("I11").Value = 0.262323813236933
$ws.Range("J11").Value = 0.262323813236933
$ws.Range("M11").Value = 6.063478666666666
$ws.Range("N11").Value = 18.190436
$ws.Range("O11").Value = 0.2107126646987263
$ws.Range("P11").Value = 0.2107126646987263
$ws.Range("Q11").Value = 18.40771873686044
$ws.Range("R11").Value = 165.669468631744
$ws.Range("S11").Value = 0.05527494970108517
$ws.Range("T11").Value = 0.05527494970108517

$ws.Range("I12").Value = 0.262323813236933
$ws.Range("J12").Value = 0.262323813236933
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.5298106666666667
$ws.Range("N12").Value = 1.589432
$ws.Range("O12").Value = 0.01841151317524362
$ws.Range("P12").Value = 0.01841151317524363
$ws.Range("Q12").Value = 1.608417588636444
$ws.Range("R12").Value = 14.475758297728
$ws.Range("S12").Value = 0.00482977834359194
$ws.Range("T12").Value = 0.004829778343591941

$ws.Range("I13").Value = 0.262323813236933
$ws.Range("J13").Value = 0.262323813236933
$ws.Range("M13").Value = 2.372966333333333
$ws.Range("N13").Value = 7.118898999999999
$ws.Range("O13").Value = 0.08246323386702208
$ws.Range("P13").Value = 0.08246323386702209
$ws.Range("Q13").Value = 7.203933457566221
$ws.Range("R13").Value = 64.83540111809599
$ws.Range("S13").Value = 0.02163206995984623
$ws.Range("T13").Value = 0.02163206995984623

$ws.Range("G14").Value = 0.9167423333333334
$ws.Range("H14").Value = 2.750227
$ws.Range("I14").Value = 0.07921490167966665
$ws.Range("J14").Value = 0.07921490167966663
$ws.Range("M14").Value = 10.055569
$ws.Range("N14").Value = 30.166707
$ws.Range("O14").Value = 0.349442268297237
$ws.Range("P14").Value = 0.3494422682972371
$ws.Range("Q14").Value = 9.218365788054335
$ws.Range("R14").Value = 82.96529209248901
$ws.Range("S14").Value = 0.02768103492588532
$ws.Range("T14").Value = 0.02768103492588532

$ws.Range("G15").Value = 0.9167423333333334
$ws.Range("H15").Value = 2.750227
$ws.Range("I15").Value = 0.07921490167966665
$ws.Range("J15").Value = 0.07921490167966663
$ws.Range("O15").Value = 0.2478120087748427
$ws.Range("P15").Value = 0.2478120087748427
$ws.Range("Q15").Value = 6.537336638439778
$ws.Range("R15").Value = 58.836029745958
$ws.Range("S15").Value = 0.01963040391013985
$ws.Range("T15").Value = 0.01963040391013985

$ws.Range("G16").Value = 0.9167423333333334
$ws.Range("H16").Value = 2.750227
$ws.Range("I16").Value = 0.07921490167966665
$ws.Range("J16").Value = 0.07921490167966663
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 2.623176333333333
$ws.Range("N16").Value = 7.869529
$ws.Range("O16").Value = 0.09115831118692827
$ws.Range("P16").Value = 0.09115831118692828
$ws.Range("Q16").Value = 2.404776792564778
$ws.Range("R16").Value = 21.642991133083
$ws.Range("S16").Value = 0.007221096657956979
$ws.Range("T16").Value = 0.007221096657956979

$ws.Range("G17").Value = 0.9167423333333334
$ws.Range("H17").Value = 2.750227
$ws.Range("I17").Value = 0.07921490167966665
$ws.Range("J17").Value = 0.07921490167966663
$ws.Range("M17").Value = 6.063478666666666
$ws.Range("N17").Value = 18.190436
$ws.Range("O17").Value = 0.2107126646987263
$ws.Range("P17").Value = 0.2107126646987263
$ws.Range("Q17").Value = 5.558647580996888
$ws.Range("R17").Value = 50.027828228972
$ws.Range("S17").Value = 0.01669158301677016
$ws.Range("T17").Value = 0.01669158301677016

$ws.Range("G18").Value = 0.9167423333333334
$ws.Range("H18").Value = 2.750227
$ws.Range("I18").Value = 0.07921490167966665
$ws.Range("J18").Value = 0.07921490167966663
$ws.Range("K18").Value = 1
$ws.Range("L18").Value = 0.3333333333333333
$ws.Range("M18").Value = 0.5298106666666667
$ws.Range("N18").Value = 1.589432
$ws.Range("O18").Value = 0.01841151317524362
$ws.Range("P18").Value = 0.01841151317524363
$ws.Range("Q18").Value = 0.4856998667848889
$ws.Range("R18").Value = 4.371298801064
$ws.Range("S18").Value = 0.001458466205950811
$ws.Range("T18").Value = 0.001458466205950811

$ws.Range("G19").Value = 0.9167423333333334
$ws.Range("H19").Value = 2.750227
$ws.Range("I19").Value = 0.07921490167966665
$ws.Range("J19").Value = 0.07921490167966663
$ws.Range("M19").Value = 2.372966333333333
$ws.Range("N19").Value = 7.118898999999999
$ws.Range("O19").Value = 0.08246323386702208
$ws.Range("P19").Value = 0.08246323386702209
$ws.Range("Q19").Value = 2.175398693341444
$ws.Range("R19").Value = 19.578588240073
$ws.Range("S19").Value = 0.006532316962963511
$ws.Range("T19").Value = 0.006532316962963511
